$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $text) {
    $range = $ws.Range($cellAddress)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue "D2" '27.164.34'
Set-TextValue "E2" '  +1.32%  '
Set-TextValue "D3" '1.645.84'
Set-TextValue "E3" '  +0.18%  '
Set-TextValue "E4" '  -0.14%  '
Set-TextValue "D5" '217.11'
Set-TextValue "E5" '  +0.09%  '
Set-TextValue "E6" '  +2.25%  '
Set-TextValue "E8" '  +1.29%  '
Set-TextValue "D9" '0.0627'
Set-TextValue "E9" '  +1.27%  '
Set-TextValue "D10" '19.91'
Set-TextValue "E10" '  +1.23%  '
Set-TextValue "D11" '0.0849'
Set-TextValue "E11" '  +0.57%  '
Set-TextValue "D12" '1.878.47'
Set-TextValue "D13" '1.653.21'
Set-TextValue "E13" '  -0.03%  '
Set-TextValue "D14" '4.14'
Set-TextValue "E14" '  +0.44%  '
Set-TextValue "D15" '0.540'
Set-TextValue "E15" '  +2.57%  '
Set-TextValue "D16" '67.56'
Set-TextValue "E16" '  +2.19%  '
Set-TextValue "D17" '27.171.03'
Set-TextValue "E17" '  +1.20%  '
Set-TextValue "D18" '0.0₃0738'
Set-TextValue "E18" '  +1.23%  '
Set-TextValue "D19" '218.68'
Set-TextValue "E19" '  +0.51%  '
Set-TextValue "E20" '  -0.16%  '
Set-TextValue "E21" '  +5.67%  '
Set-TextValue "D22" '6.82'
Set-TextValue "E22" '  +2.83%  '
Set-TextValue "D23" '4.40'
Set-TextValue "E23" '  +0.33%  '
Set-TextValue "D24" '9.20'
Set-TextValue "E24" '  +0.45%  '
Set-TextValue "D25" '147.80'
Set-TextValue "E25" '  +1.26%  '
Set-TextValue "D26" '7.56'
Set-TextValue "E26" '  +3.14%  '
Set-TextValue "E27" '  -0.04%  '
Set-TextValue "E28" '  -0.05%  '
Set-TextValue "D29" '15.74'
Set-TextValue "E29" '  -0.49%  '
Set-TextValue "D30" '0.0507'
Set-TextValue "E30" '  -0.62%  '
Set-TextValue "E31" '  +0.14%  '
Set-TextValue "D32" '3.39'
Set-TextValue "E32" '  +0.57%  '
Set-TextValue "D33" '3.03'
Set-TextValue "E33" '  +1.32%  '
Set-TextValue "E34" '  +1.53%  '
Set-TextValue "D35" '1.263.00'
Set-TextValue "E35" '  +1.42%  '
Set-TextValue "E36" '  +0.41%  '
Set-TextValue "D37" '0.0177'
Set-TextValue "E37" '  +1.56%  '
Set-TextValue "D38" '0.548'
Set-TextValue "E38" '  +2.54%  '
Set-TextValue "D39" '0.848'
Set-TextValue "E39" '  +2.03%  '
Set-TextValue "E40" '  -0.11%  '
Set-TextValue "E41" '  +0.25%  '
Set-TextValue "E42" '  +6.08%  '
Set-TextValue "D43" '5.43'
Set-TextValue "E43" '  +1.54%  '
Set-TextValue "D44" '1.788.09'
Set-TextValue "E44" '  +0.13%  '
Set-TextValue "D45" '61.87'
Set-TextValue "E45" '  +1.67%  '
Set-TextValue "D46" '91.44'
Set-TextValue "E46" '  -0.17%  '
Set-TextValue "E47" '  +0.84%  '
Set-TextValue "D48" '0.0₆0108'
Set-TextValue "E48" '  +1.67%  '
Set-TextValue "D49" '0.0514'
Set-TextValue "E49" '  +0.02%  '
Set-TextValue "D50" '7.67'
Set-TextValue "E50" '  +2.10%  '
Set-TextValue "D51" '0.0972'
Set-TextValue "E51" '  +0.12%  '
